$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.010.80"
$ws.Range("E2").Value = "  +2.50%  "
$ws.Range("D3").Value = "3.574.47"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  +1.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "190.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.627"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.37%  "
$ws.Range("D8").Value = "3.568.80"
$ws.Range("E8").Value = "  +2.02%  "
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.219"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +15.90%  "
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.90"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000320"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.13%  "
$ws.Range("E14").Value = "  +1.26%  "
$ws.Range("D15").Value = "4.144.05"
$ws.Range("E15").Value = "  +1.92%  "
$ws.Range("D16").Value = "71.022.86"
$ws.Range("E16").Value = "  +2.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "12.85"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("D19").Value = "3.568.99"
$ws.Range("E19").Value = "  +1.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "567.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.43%  "
$ws.Range("E21").Value = "  +0.67%  "
$ws.Range("E22").Value = "  -0.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.78"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.62"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.35%  "
$ws.Range("E25").Value = "  -1.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "94.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.50%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.34"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.62"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.47%  "
$ws.Range("E31").Value = "  -0.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.90%  "
$ws.Range("E33").Value = "  +2.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "64.08"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.81"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +24.71%  "
$ws.Range("E36").Value = "  +6.71%  "
$ws.Range("B37").Value = "TheGraph"
$ws.Range("C37").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.413"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.81%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "532.87"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "38.57"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.65%  "
$ws.Range("E40").Value = "  +4.85%  "
$ws.Range("D41").Value = "3.653.52"
$ws.Range("E41").Value = "  +10.37%  "
$ws.Range("E43").Value = "  +4.66%  "
$ws.Range("E44").Value = "  +2.74%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0474"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.73%  "
$ws.Range("E46").Value = "  -0.79%  "
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.42%  "
$ws.Range("E49").Value = "  +3.47%  "
$ws.Range("E51").Value = "  +7.76%  "
